# Generate Report for Handoff
# Updates the localization-status report: the file previously tracked as
# "8eb115b1-3ae2-4155-9391-035f78d714af" is renamed/regenerated as
# "6b984f79-5315-48d5-bbca-213b847ccf8e" and is now "Ready for handoff"
# (rather than "Handed back"), and the second tracked file's id changes
# from "eb62c29d-6ef4-4e89-b977-4723e27828fb" to
# "ffff56229551-b32e-4158-8590-ca88962b40b6". The per-language sheets lose
# their "Latest Target File" / "Latest Handback File" columns (F/G) for
# both data rows, since there is no handback yet.

$wb = $excel.ActiveWorkbook

function Remove-HyperlinksAt($ws, $addrs) {
    $changed = $true
    while ($changed) {
        $changed = $false
        foreach ($h in @($ws.Hyperlinks)) {
            $a = $h.Range.Address()
            if ($addrs -contains $a) {
                $h.Delete()
                $changed = $true
                break
            }
        }
    }
}

function Set-HyperlinkDisplay($ws, $addr, $newText) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $newText
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "6b984f79-5315-48d5-bbca-213b847ccf8e.md"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-52-11 10:52:40"

$wsOverview.Range("A3").Value = "ffff56229551-b32e-4158-8590-ca88962b40b6.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-52-11 10:52:40"

Set-HyperlinkDisplay $wsOverview '$A$2' "6b984f79-5315-48d5-bbca-213b847ccf8e.md"
Set-HyperlinkDisplay $wsOverview '$A$3' "ffff56229551-b32e-4158-8590-ca88962b40b6.md"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Remove-HyperlinksAt $wsZhCn @('$F$2', '$G$2', '$F$3', '$G$3')

$wsZhCn.Range("A2").Value = "6b984f79-5315-48d5-bbca-213b847ccf8e.md"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("D2").Value = "6b984f79-5315-48d5-bbca-213b847ccf8e.7d306ebc8b8eab999ead8f0a675bfed8bcaa2154.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-11 10:52:36"
$wsZhCn.Range("H2").Value = "0001-01-01 00:00:00"

$wsZhCn.Range("A3").Value = "ffff56229551-b32e-4158-8590-ca88962b40b6.md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "6b984f79-5315-48d5-bbca-213b847ccf8e.7d306ebc8b8eab999ead8f0a675bfed8bcaa2154.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-11 10:52:36"
$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"

$wsZhCn.Range("F2:G3").Clear()

Set-HyperlinkDisplay $wsZhCn '$A$2' "6b984f79-5315-48d5-bbca-213b847ccf8e.md"
Set-HyperlinkDisplay $wsZhCn '$D$2' "6b984f79-5315-48d5-bbca-213b847ccf8e.7d306ebc8b8eab999ead8f0a675bfed8bcaa2154.zh-cn.xlf"
Set-HyperlinkDisplay $wsZhCn '$A$3' "ffff56229551-b32e-4158-8590-ca88962b40b6.md"
Set-HyperlinkDisplay $wsZhCn '$D$3' "6b984f79-5315-48d5-bbca-213b847ccf8e.7d306ebc8b8eab999ead8f0a675bfed8bcaa2154.zh-cn.xlf"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Remove-HyperlinksAt $wsDeDe @('$F$2', '$G$2', '$F$3', '$G$3')

$wsDeDe.Range("A2").Value = "6b984f79-5315-48d5-bbca-213b847ccf8e.md"
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("D2").Value = "6b984f79-5315-48d5-bbca-213b847ccf8e.7d306ebc8b8eab999ead8f0a675bfed8bcaa2154.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-11 10:52:40"
$wsDeDe.Range("H2").Value = "0001-01-01 00:00:00"

$wsDeDe.Range("A3").Value = "ffff56229551-b32e-4158-8590-ca88962b40b6.md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "6b984f79-5315-48d5-bbca-213b847ccf8e.7d306ebc8b8eab999ead8f0a675bfed8bcaa2154.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-11 10:52:40"
$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"

$wsDeDe.Range("F2:G3").Clear()

Set-HyperlinkDisplay $wsDeDe '$A$2' "6b984f79-5315-48d5-bbca-213b847ccf8e.md"
Set-HyperlinkDisplay $wsDeDe '$D$2' "6b984f79-5315-48d5-bbca-213b847ccf8e.7d306ebc8b8eab999ead8f0a675bfed8bcaa2154.de-de.xlf"
Set-HyperlinkDisplay $wsDeDe '$A$3' "ffff56229551-b32e-4158-8590-ca88962b40b6.md"
Set-HyperlinkDisplay $wsDeDe '$D$3' "6b984f79-5315-48d5-bbca-213b847ccf8e.7d306ebc8b8eab999ead8f0a675bfed8bcaa2154.de-de.xlf"
